$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MapSet")

$row = 23

# Set cells in an order that matches the shared-string insertion order
# expected by the target file: qiongqi (218), 穷奇 (219), Qiongqi (220)
$ws.Cells.Item($row, 1).Value = 42050003           # A Id
$ws.Cells.Item($row, 6).Value = "qiongqi"          # F Ename  -> new shared string "qiongqi"
$ws.Cells.Item($row, 7).Value = "qiongqi"          # G Figue
$ws.Cells.Item($row, 8).Value = "qiongqi"          # H Script
$ws.Cells.Item($row, 2).Value = "穷奇"              # B Name   -> new shared string "穷奇"
$ws.Cells.Item($row, 21).Value = "Qiongqi"         # U Flag   -> new shared string "Qiongqi"
$ws.Cells.Item($row, 3).Value = 2                  # C Type
$ws.Cells.Item($row, 4).Value = 0                  # D Level
$ws.Cells.Item($row, 5).Value = 3                  # E Danger
$ws.Cells.Item($row, 9).Value = "true"             # I TriggerMulti
$ws.Cells.Item($row, 10).Value = 43000019          # J EnemyId
$ws.Cells.Item($row, 11).Value = "oneline"         # K BattleMap
$ws.Cells.Item($row, 12).Value = 13000002          # L SceneId
$ws.Cells.Item($row, 13).Value = 300               # M RewardGold
$ws.Cells.Item($row, 14).Value = 500               # N RewardFood
$ws.Cells.Item($row, 15).Value = 500               # O RewardHealth
$ws.Cells.Item($row, 16).Value = 500               # P RewardMental
$ws.Cells.Item($row, 17).Value = 300               # Q RewardExp
$ws.Cells.Item($row, 26).Value = 200               # Z PunishHealth
$ws.Cells.Item($row, 27).Value = 200               # AA PunishMental

# copy styles from row 22 to keep formatting consistent with existing rows
$ws.Range("A22:L22").Copy() | Out-Null
$ws.Range("A23:L23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A1").Select() | Out-Null

$sel = $ws.Range("A23")
$sel.Select() | Out-Null

# Update table range to include new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:AF23"))
